$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: swap text so each column now shows its counterpart label
$ws.Range("B1").Value = "FFR_CA"
$ws.Range("C1").Value = "LF_CA"
$ws.Range("D1").Value = "CA_FFR"
$ws.Range("E1").Value = "CA_LF"

# Update data row 2 values
$ws.Range("B2").Value = 0.05716633442687229
$ws.Range("C2").Value = 0.4868474009768662
$ws.Range("D2").Value = 14.02068439442374
$ws.Range("E2").Value = 0.3511567374768891
